# Generate Report for Handoff
#
# Regenerates the handoff/handback report timestamps and priority values:
#  - "Priority" column (E) for the md/ht rows that were still "low" is now "ht"
#    on both the zh-cn and de-de sheets.
#  - The zh-cn "Latest Handoff Datetime" (H) for those same rows advances
#    from 04:30:23 to 04:30:38.
#  - The "Latest HO Xliff Generate Date" (Overview!G and de-de!H for those
#    same rows) advances from 04:30:28 to 04:30:43.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Rows 4-7 on the zh-cn / de-de sheets correspond to the 3a2c05e0, ad56e393,
# b0f27474 and c25729fe entries, whose Priority is still "low".
$rows = 4, 5, 6, 7

foreach ($r in $rows) {
    $wsZhCn.Range("E$r").Value = "ht"
    $wsDeDe.Range("E$r").Value = "ht"

    $wsZhCn.Range("H$r").Value = "2016-08-25 04:30:38"

    $wsOverview.Range("G$r").Value = "2016-08-25 04:30:43"
    $wsDeDe.Range("H$r").Value = "2016-08-25 04:30:43"
}
